$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

# Each data row in the table holds 5 problems (columns 1-5); the data rows
# are rows 1, 5, 9, 13, 17 (the rows in between are blank spacer rows).
$rows = @(1, 5, 9, 13, 17)

$values = @(
    @("50÷3=16, 2", "48÷5=9, 3", "87÷8=10, 7", "31÷2=15, 1", "69÷4=17, 1"),
    @("62÷9=6, 8", "45÷8=5, 5", "33÷5=6, 3", "25÷5=5, 0", "22÷7=3, 1"),
    @("37÷5=7, 2", "51÷2=25, 1", "49÷5=9, 4", "84÷8=10, 4", "84÷9=9, 3"),
    @("38÷5=7, 3", "66÷2=33, 0", "30÷6=5, 0", "11÷2=5, 1", "80÷4=20, 0"),
    @("20÷9=2, 2", "21÷6=3, 3", "34÷4=8, 2", "29÷7=4, 1", "79÷8=9, 7")
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $rows[$i]
    for ($c = 1; $c -le 5; $c++) {
        $cell = $tbl.Cell($r, $c)
        $cell.Range.Text = $values[$i][$c - 1]
    }
}
